$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The May-9th data collection added 7 new samples at the start of the series
# and 3 new samples at the end. Shift the existing 20 data rows (2-21) down by 7
# rows (landing on 9-28) working bottom-up so we never overwrite an unread source row.
for ($r = 21; $r -ge 2; $r--) {
    $dest = $r + 7
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value()
}

# Fill the newly-opened rows 2-8 with the new samples that now lead the series
$ws.Cells.Item(2, 1).Value = 1.327703475952148
$ws.Cells.Item(2, 2).Value = -2.356002759933471
$ws.Cells.Item(2, 3).Value = -1.412894463539124
$ws.Cells.Item(3, 1).Value = 1.294596195220947
$ws.Cells.Item(3, 2).Value = -2.42879855632782
$ws.Cells.Item(3, 3).Value = -1.538487493991852
$ws.Cells.Item(4, 1).Value = 1.5967435836792
$ws.Cells.Item(4, 2).Value = -2.721291160583497
$ws.Cells.Item(4, 3).Value = -1.309774732589722
$ws.Cells.Item(5, 1).Value = 1.538975667953491
$ws.Cells.Item(5, 2).Value = -3.117227482795715
$ws.Cells.Item(5, 3).Value = -1.760656356811524
$ws.Cells.Item(6, 1).Value = 1.350284004211427
$ws.Cells.Item(6, 2).Value = -3.035107040405273
$ws.Cells.Item(6, 3).Value = -2.307204818725585
$ws.Cells.Item(7, 1).Value = 2.076921081542969
$ws.Cells.Item(7, 2).Value = -3.024871301651001
$ws.Cells.Item(7, 3).Value = -2.036388444900513
$ws.Cells.Item(8, 1).Value = 2.578349113464355
$ws.Cells.Item(8, 2).Value = -2.735702991485596
$ws.Cells.Item(8, 3).Value = -2.153444766998291

# Append the 3 new samples collected at the end of the series (rows 29-31)
$ws.Cells.Item(29, 1).Value = 2.200880432128908
$ws.Cells.Item(29, 2).Value = -3.761867809295656
$ws.Cells.Item(29, 3).Value = -0.7659695267677282
$ws.Cells.Item(30, 1).Value = 2.088931465148925
$ws.Cells.Item(30, 2).Value = -3.697214221954345
$ws.Cells.Item(30, 3).Value = -0.6586695432662961
$ws.Cells.Item(31, 1).Value = 1.876247692108153
$ws.Cells.Item(31, 2).Value = -2.620123100280754
$ws.Cells.Item(31, 3).Value = -1.525603616237647
